# Populate column B with the per-hospital bed/unit counts next to the
# existing hospital names in column A, then add a total row with a SUM
# formula, and finally update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$counts = @{
    2  = 23
    3  = 3
    4  = 45
    5  = 23
    6  = 12
    7  = 3
    8  = 12
    9  = 4
    10 = 32
    11 = 21
    12 = 1
    13 = 23
    14 = 11
    15 = 32
    16 = 12
    17 = 21
    18 = 15
    19 = 8
    20 = 9
    21 = 4
    22 = 7
    23 = 14
    24 = 5
    25 = 10
    26 = 11
    27 = 4
    28 = 12
    29 = 15
    30 = 12
    31 = 23
    32 = 12
    33 = 13
    34 = 14
    35 = 5
    36 = 12
    37 = 10
    38 = 7
}

foreach ($row in $counts.Keys) {
    $ws.Cells.Item($row, 2).Value = $counts[$row]
}

# Total row with a SUM formula over the new data.
$ws.Range("B39").Formula = "=SUM(B2:B38)"

# Match the recorded selection after the edit.
$ws.Range("C40").Select()
